# Calibration of energy use modeling by renovation level
#
# Sheet1 keeps its original raw survey shares (% SFH / % TH / % AB).
# Sheet2 is a calibration pass: % TH is rescaled (*2/3*0.5) and % SFH is
#   backed out as the remainder (1 - %TH - %AB), with a 12% bump applied
#   to % AB for the two newest age classes.
# Sheet3 is that calibrated table pasted down to plain values.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Build Sheet2 (formula-driven calibration) as a copy of Sheet1 ---
$ws1.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Name = "Sheet2"

$ws2.Range("C2").Formula = "=0.02*2/3*0.5"
$ws2.Range("C3").Formula = "=0.197740440564393*2/3*0.5"
$ws2.Range("C4").Formula = "=0.189333525195886*2/3*0.5"
$ws2.Range("C5").Formula = "=0.189063066210266*2/3*0.5"
$ws2.Range("C6").Formula = "=0.213790952960481*2/3*0.5"
$ws2.Range("C7").Formula = "=0.241475102437669*2/3*0.5"
$ws2.Range("C8").Formula = "=0.23047848210847*2/3*0.5"
$ws2.Range("C9").Formula = "=0.247551850996143*2/3*0.5"

$ws2.Range("D8").Formula = "=0.432760636249306*1.12"
$ws2.Range("D9").Formula = "=0.412213197789113*1.12"

$ws2.Range("B2").Formula = "=1-C2"
$ws2.Range("B3").Formula = "=1-C3-D3"
$ws2.Range("B4:B9").Formula = "=1-C4-D4"

$ws2.Range("A1:D9").Select() | Out-Null

# --- Build Sheet3 (values-only snapshot of the calibrated table) ---
$ws1.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Sheet3"

# (Literal decimals below are the same numbers as Sheet2's computed results,
#  written without E-notation since the COM parser only accepts plain decimals.)
$ws3.Range("B2").Value = 0.99333333333333329040
$ws3.Range("C2").Value = 0.00666666666666666709
$ws3.Range("D2").Value = 0

$ws3.Range("B3").Value = 0.68606232286796176023
$ws3.Range("C3").Value = 0.06591348018813099341
$ws3.Range("D3").Value = 0.24802419694390720473

$ws3.Range("B4").Value = 0.65522686199512292227
$ws3.Range("C4").Value = 0.06311117506529533017
$ws3.Range("D4").Value = 0.28166196293958178920

$ws3.Range("B5").Value = 0.69873203660591454156
$ws3.Range("C5").Value = 0.06302102207008866697
$ws3.Range("D5").Value = 0.23824694132399679147

$ws3.Range("B6").Value = 0.75474946364816108968
$ws3.Range("C6").Value = 0.07126365098682700683
$ws3.Range("D6").Value = 0.17398688536501188961

$ws3.Range("B7").Value = 0.63712757830404898041
$ws3.Range("C7").Value = 0.08049170081255634102
$ws3.Range("D7").Value = 0.28238072088339472021

$ws3.Range("B8").Value = 0.43848192669795388321
$ws3.Range("C8").Value = 0.07682616070282333176
$ws3.Range("D8").Value = 0.48469191259922278503

$ws3.Range("B9").Value = 0.45580393481081238338
$ws3.Range("C9").Value = 0.08251728366538099713
$ws3.Range("D9").Value = 0.46167878152380659174

$ws3.Range("B2:D9").Select() | Out-Null

# --- Sheet1: drop the stale single-cell selection, select the full used range instead ---
$ws1.Rows("1:9").Select() | Out-Null

# Sheet3 ends up the active/visible tab, matching the authored workbook state.
$ws3.Activate() | Out-Null
